$wb = $excel.ActiveWorkbook

# --- Update selections on existing sheets (Slovakia, Spain) ---
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Activate()
$slovakia.Cells.Select()

$spain = $wb.Worksheets.Item("Spain")
$spain.Activate()
$spain.Range("C13").Select()

# --- Add the new "Croatia" sheet by duplicating "Slovakia" (same layout/styles) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $lastSheet)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# --- Fill in the Croatia-specific content ---
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2416"

# --- Make Croatia the active/selected sheet with B4 selected ---
$croatia.Activate()
$croatia.Range("B4").Select()
